$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "Formula" column
$ws.Range("H1").Value = "Formula"

# Fill in the formula column: H = G + A, for each data row (2..46).
# H2 is entered on its own (matches the author's first keystroke), then
# H3:H46 is filled as a second pass, mirroring how Excel records the
# resulting shared-formula group (master cell H3, range H3:H46).
$ws.Range("H2").Formula = "=G2+A2"
$ws.Range("H3:H46").Formula = "=G3+A3"

# Update the view: scroll so row 6 is the top-left visible row, and select L35
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L35").Select()
